$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.952.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.94%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.963.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.20%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'595.35"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'146.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.26%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.10%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'2.962.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.21%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +0.81%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +3.06%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +7.05%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.80%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +6.62%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'33.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.09%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = "'  -0.51%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.454.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.16%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'62.806.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.91%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "'Polkadot"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'6.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.15%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = "'WrappedEther"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'2.948.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.44%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'442.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.26%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'13.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.50%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.672"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.15%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.01%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'11.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.43%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'81.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.35%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.21%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +0.17%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +0.04%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'7.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.30%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  +0.52%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  -3.11%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.0₃0973"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +10.57%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -0.68%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'26.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.73%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +0.08%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.55%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'Filecoin"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'5.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.49%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'dogwifhat"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'3.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +5.09%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.95%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'49.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.94%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'8.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.12%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.118"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.26%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'40.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.16%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +0.27%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.742.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.52%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'134.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.62%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'366.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.47%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.0340"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.22%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +0.05%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'Stellar"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.23%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'InjectiveProtocol"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'23.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.52%  "
$ws.Range("E51").Style = "Normal"

